$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 39583.668
$ws.Range("I64").Value = 45500.4
$ws.Range("K64").Value = 45500.4
$ws.Range("M64").Value = -45252.4
$ws.Range("H67").Value = 39583.668
$ws.Range("I67").Value = 45500.4
$ws.Range("K67").Value = 45500.4
$ws.Range("M67").Value = -44642.4
$ws.Range("H94").Value = 1922.5
$ws.Range("I94").Value = 1922.5
$ws.Range("K94").Value = 1922.5
$ws.Range("M94").Value = -1471.5
$ws.Range("H98").Value = 4252.0713
$ws.Range("I98").Value = 3040.2307
$ws.Range("K98").Value = 3040.2307
$ws.Range("M98").Value = -1542.2307
$ws.Range("H100").Value = 2311.0667
$ws.Range("J100").Value = 800
$ws.Range("L100").Value = 800
$ws.Range("N100").Value = -1882
$ws.Range("H113").Value = 6332.375
$ws.Range("I113").Value = 6307.7856
$ws.Range("K113").Value = 6307.7856
$ws.Range("M113").Value = -3053.7856
$ws.Range("H122").Value = 4252.0713
$ws.Range("I122").Value = 3040.2307
$ws.Range("K122").Value = 9120.6921
$ws.Range("M122").Value = -6670.6921
$ws.Range("H132").Value = 3264087
$ws.Range("I132").Value = 3502834.8
$ws.Range("K132").Value = 10508504.4
$ws.Range("M132").Value = -10505974.4

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20576.4
$ws.Range("I32").Value = 21659.785
$ws.Range("K32").Value = 21659.785
$ws.Range("M32").Value = -21372.785
$ws.Range("H43").Value = 30278.715
$ws.Range("I43").Value = 20342
$ws.Range("K43").Value = 20342
$ws.Range("M43").Value = -20029
$ws.Range("H52").Value = 50000
$ws.Range("I52").Value = 50000
$ws.Range("K52").Value = 50000
$ws.Range("M52").Value = -49682
$ws.Range("H61").Value = 21596.4
$ws.Range("I61").Value = 994
$ws.Range("K61").Value = 994
$ws.Range("M61").Value = -782
$ws.Range("H63").Value = 4161.636
$ws.Range("I63").Value = 4130
$ws.Range("K63").Value = 4130
$ws.Range("M63").Value = -3444
$ws.Range("H66").Value = 4161.636
$ws.Range("I66").Value = 4130
$ws.Range("K66").Value = 20650
$ws.Range("M66").Value = -17218
$ws.Range("H74").Value = 510660.1
$ws.Range("I74").Value = 1200901.4
$ws.Range("K74").Value = 1200901.4
$ws.Range("M74").Value = -1200027.4
$ws.Range("H77").Value = 510660.1
$ws.Range("I77").Value = 1200901.4
$ws.Range("K77").Value = 6004507
$ws.Range("M77").Value = -6000139
$ws.Range("H97").Value = 727.06665
$ws.Range("I97").Value = 640.9524
$ws.Range("K97").Value = 640.9524
$ws.Range("M97").Value = -144.9524
$ws.Range("H102").Value = 2999.25
$ws.Range("I102").Value = 2999.5
$ws.Range("J102").Value = 2999
$ws.Range("K102").Value = 2999.5
$ws.Range("L102").Value = 2999
$ws.Range("M102").Value = -1377.5
$ws.Range("N102").Value = -6243
$ws.Range("H132").Value = 3157.5557
$ws.Range("I132").Value = 1974
$ws.Range("K132").Value = 5922
$ws.Range("M132").Value = -3392
$ws.Range("H136").Value = 21596.4
$ws.Range("I136").Value = 994
$ws.Range("K136").Value = 2982
$ws.Range("M136").Value = -432

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14761.375
$ws.Range("I82").Value = 11870.143
$ws.Range("K82").Value = 11870.143
$ws.Range("M82").Value = -11487.143
$ws.Range("H85").Value = 14761.375
$ws.Range("I85").Value = 11870.143
$ws.Range("K85").Value = 11870.143
$ws.Range("M85").Value = -10544.143
$ws.Range("H94").Value = 5373.55
$ws.Range("I94").Value = 6595.9375
$ws.Range("J94").Value = 484
$ws.Range("K94").Value = 6595.9375
$ws.Range("L94").Value = 484
$ws.Range("M94").Value = -6144.9375
$ws.Range("N94").Value = -1386
$ws.Range("H105").Value = 2260.5
$ws.Range("I105").Value = 2121.0476
$ws.Range("J105").Value = 3236.6667
$ws.Range("K105").Value = 2121.0476
$ws.Range("L105").Value = 3236.6667
$ws.Range("M105").Value = -374.0475999999999
$ws.Range("N105").Value = -6730.6667
$ws.Range("H134").Value = 1572.1765
$ws.Range("I134").Value = 1357.9375
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 4073.8125
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -1538.8125
$ws.Range("N134").Value = -20070

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H22").Value = 681.88464
$ws.Range("I22").Value = 605.75
$ws.Range("J22").Value = 747.1429000000001
$ws.Range("K22").Value = 605.75
$ws.Range("L22").Value = 747.1429000000001
$ws.Range("M22").Value = -255.75
$ws.Range("N22").Value = -1447.1429
$ws.Range("H31").Value = 7143797
$ws.Range("I31").Value = 7143797
$ws.Range("K31").Value = 7143797
$ws.Range("M31").Value = -7143502
$ws.Range("H34").Value = 7143797
$ws.Range("I34").Value = 7143797
$ws.Range("K34").Value = 7143797
$ws.Range("M34").Value = -7143595
$ws.Range("H58").Value = 1984.4762
$ws.Range("I58").Value = 1118.0667
$ws.Range("J58").Value = 4150.5
$ws.Range("K58").Value = 1118.0667
$ws.Range("L58").Value = 4150.5
$ws.Range("M58").Value = -915.0667000000001
$ws.Range("N58").Value = -4556.5
$ws.Range("H132").Value = 44563.125
$ws.Range("I132").Value = 79285.766
$ws.Range("J132").Value = 3527.2727
$ws.Range("K132").Value = 237857.298
$ws.Range("L132").Value = 10581.8181
$ws.Range("M132").Value = -235327.298
$ws.Range("N132").Value = -15641.8181
$ws.Range("H136").Value = 1984.4762
$ws.Range("I136").Value = 1118.0667
$ws.Range("J136").Value = 4150.5
$ws.Range("K136").Value = 3354.2001
$ws.Range("L136").Value = 12451.5
$ws.Range("M136").Value = -804.2001
$ws.Range("N136").Value = -17551.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 127569810
$ws.Range("I4").Value = 71095976
$ws.Range("K4").Value = 213287928
$ws.Range("M4").Value = -213287816
$ws.Range("H68").Value = 4894.684
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 4894.684
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 14684.052
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -16306.052
$ws.Range("H71").Value = 4894.684
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4894.684
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 44052.156
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -52164.156
$ws.Range("H92").Value = 1438.9286
$ws.Range("J92").Value = 1516.3334
$ws.Range("L92").Value = 4549.0002
$ws.Range("N92").Value = -7045.0002
$ws.Range("H97").Value = 873
$ws.Range("I97").Value = 722.5
$ws.Range("J97").Value = 910.625
$ws.Range("K97").Value = 2167.5
$ws.Range("L97").Value = 2731.875
$ws.Range("M97").Value = -1671.5
$ws.Range("N97").Value = -3723.875
$ws.Range("H107").Value = 1229.3462
$ws.Range("I107").Value = 1524.1333
$ws.Range("J107").Value = 827.36365
$ws.Range("K107").Value = 4572.3999
$ws.Range("L107").Value = 2482.09095
$ws.Range("M107").Value = -2652.3999
$ws.Range("N107").Value = -6322.09095
$ws.Range("H133").Value = 1500
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H137").Value = 2729.9285
$ws.Range("I137").Value = 2515.5
$ws.Range("J137").Value = 4016.5
$ws.Range("K137").Value = 7546.5
$ws.Range("L137").Value = 12049.5
$ws.Range("M137").Value = -2446.5
$ws.Range("N137").Value = -22249.5
$ws.Range("H138").Value = 12981
$ws.Range("I138").Value = 10963.25
$ws.Range("J138").Value = 17016.5
$ws.Range("K138").Value = 32889.75
$ws.Range("L138").Value = 51049.5
$ws.Range("M138").Value = -27749.75
$ws.Range("N138").Value = -61329.5
$ws.Range("H139").Value = 6300.2856
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 2732.182
$ws.Range("I140").Value = 2732.182
$ws.Range("K140").Value = 8196.545999999998
$ws.Range("M140").Value = -3016.545999999998

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 822.25
$ws.Range("I107").Value = 846.5
$ws.Range("K107").Value = 846.5
$ws.Range("M107").Value = 1073.5
$ws.Range("H126").Value = 4425.273
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 4964.222
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 14892.666
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -19832.666
$ws.Range("H132").Value = 2611.5
$ws.Range("I132").Value = 2302.1428
$ws.Range("J132").Value = 3333.3333
$ws.Range("K132").Value = 6906.428400000001
$ws.Range("L132").Value = 9999.999899999999
$ws.Range("M132").Value = -4376.428400000001
$ws.Range("N132").Value = -15059.9999
$ws.Range("H140").Value = 35000
$ws.Range("J140").Value = 48888.89
$ws.Range("L140").Value = 48888.89
$ws.Range("N140").Value = -59248.89

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10011667
$ws.Range("J2").Value = 10011667
$ws.Range("L2").Value = 10011667
$ws.Range("N2").Value = -10011891

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 724.5
$ws.Range("I100").Value = 734.0769
$ws.Range("J100").Value = 600
$ws.Range("K100").Value = 1468.1538
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = -927.1538
$ws.Range("N100").Value = -2282
$ws.Range("H132").Value = 352500
$ws.Range("I132").Value = 700000
$ws.Range("K132").Value = 2100000
$ws.Range("M132").Value = -2097470
$ws.Range("H136").Value = 35471.61
$ws.Range("J136").Value = 2588.8
$ws.Range("L136").Value = 7766.400000000001
$ws.Range("N136").Value = -12866.4
